# Generate Report for Handback
# Adds a new handback record (a35f5189-d9b8-4f18-a057-abf087767076.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" worksheets, expanding each
# sheet's table/autofilter/dimension from 3 to 4 data rows and wiring up
# the matching hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet (row 4)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.md"
$wsOverview.Range("B4").Value = "e2e\a35f5189-d9b8-4f18-a057-abf087767076.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-21 18:49:00"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/483fdf26001d463f7b5b949236183b6c43501f54/e2e/a35f5189-d9b8-4f18-a057-abf087767076.md", "", "", "e2e\a35f5189-d9b8-4f18-a057-abf087767076.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# zh-cn sheet (row 4)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.483fdf26001d463f7b5b949236183b6c43501f54.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-21 18:48:55"
$wsZhCn.Range("I4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.md"
$wsZhCn.Range("J4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.483fdf26001d463f7b5b949236183b6c43501f54.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-21 18:49:24"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("O4").Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/483fdf26001d463f7b5b949236183b6c43501f54/e2e/a35f5189-d9b8-4f18-a057-abf087767076.md", "", "", "a35f5189-d9b8-4f18-a057-abf087767076.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/483fdf26001d463f7b5b949236183b6c43501f54/e2e/a35f5189-d9b8-4f18-a057-abf087767076.md", "", "", "a35f5189-d9b8-4f18-a057-abf087767076.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------
# de-de sheet (row 4)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.483fdf26001d463f7b5b949236183b6c43501f54.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-21 18:49:00"
$wsDeDe.Range("I4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.md"
$wsDeDe.Range("J4").Value = "a35f5189-d9b8-4f18-a057-abf087767076.483fdf26001d463f7b5b949236183b6c43501f54.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-21 18:49:31"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("O4").Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/483fdf26001d463f7b5b949236183b6c43501f54/e2e/a35f5189-d9b8-4f18-a057-abf087767076.md", "", "", "a35f5189-d9b8-4f18-a057-abf087767076.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/483fdf26001d463f7b5b949236183b6c43501f54/e2e/a35f5189-d9b8-4f18-a057-abf087767076.md", "", "", "a35f5189-d9b8-4f18-a057-abf087767076.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))

# ---------------------------------------------------------------------
# Match existing visual styling: hyperlink font for the link cells and
# the yyyy-mm-dd HH:mm:ss date format used for the "datetime" columns.
# ---------------------------------------------------------------------
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Host "Handback row for a35f5189-d9b8-4f18-a057-abf087767076.md added."
